# Fixed bug in saving function
#
# 1) "shunt_Y_profiles" sheet: the per-snapshot Y admittance profile for the
#    single shunt ("shunt1@Bus 9 LV") was not being written out. Fill column B
#    (rows 2-11) with the shunt's Y value ("19j"), matching the "shunt" sheet.
# 2) "bus" sheet: the saved bus-diagram coordinates (x/y, columns I/J) had
#    drifted - shift them back to the correct position.
# 3) Force a full recalculation the next time the workbook is opened.

$wb = $excel.ActiveWorkbook

# --- 1) shunt_Y_profiles: stamp the Y value onto every profile row ---------
$wsShuntYProfiles = $wb.Worksheets.Item("shunt_Y_profiles")
for ($r = 2; $r -le 11; $r++) {
    $wsShuntYProfiles.Cells.Item($r, 2).Value = "19j"
}

# --- 2) bus: correct the drifted diagram coordinates ------------------------
# x/y (columns I/J) are stored as plain text (not numbers), so force the
# "@" text format before writing, otherwise Excel re-interprets the literal
# as a number and drops the trailing ".0". Reset the style back to Normal
# afterwards so the cells keep their original (default) look.
$wsBus = $wb.Worksheets.Item("bus")

$xValues = @("-2246.0", "-1993.0", "-1776.0", "-1874.0", "-2217.0", "-2891.0", "-1916.0", "-1802.0", "-2181.0", "-2450.0", "-2649.0", "-3507.0", "-3097.0", "-2701.0")
$yValues = @("-549.0", "-648.0", "-659.0", "-1016.0", "-795.0", "-722.0", "-1291.0", "-1563.0", "-1248.0", "-1192.0", "-1044.0", "-1112.0", "-1290.0", "-1446.0")

$coordRange = $wsBus.Range("I2:J15")
$coordRange.NumberFormat = "@"

for ($i = 0; $i -lt $xValues.Length; $i++) {
    $row = $i + 2
    $newX = [double]$xValues[$i] - 856.0
    $newY = [double]$yValues[$i] - 520.0
    $wsBus.Cells.Item($row, 9).Value = $newX.ToString("0.0")
    $wsBus.Cells.Item($row, 10).Value = $newY.ToString("0.0")
}

$coordRange.Style = "Normal"

# --- 3) force full recalculation on next load --------------------------------
$wb.ForceFullCalculation = $true
